# Update "想去人数" (interested-count) figures in column F for the
# "展览" and "全部类型" sheets, matching the latest scrape snapshot.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 14012
$ws1.Range("F6").Value = 1808
$ws1.Range("F7").Value = 173
$ws1.Range("F8").Value = 102
$ws1.Range("F10").Value = 563
$ws1.Range("F11").Value = 38
$ws1.Range("F14").Value = 14240
$ws1.Range("F15").Value = 382
$ws1.Range("F16").Value = 640
$ws1.Range("F17").Value = 15060
$ws1.Range("F19").Value = 8416
$ws1.Range("F20").Value = 293
$ws1.Range("F21").Value = 4
$ws1.Range("F23").Value = 165
$ws1.Range("F38").Value = 15
$ws1.Range("F43").Value = 5185

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 14012
$ws4.Range("F6").Value = 1808
$ws4.Range("F7").Value = 173
$ws4.Range("F8").Value = 102
$ws4.Range("F10").Value = 563
$ws4.Range("F11").Value = 38
$ws4.Range("F14").Value = 14240
$ws4.Range("F15").Value = 382
$ws4.Range("F16").Value = 640
$ws4.Range("F17").Value = 15061
$ws4.Range("F19").Value = 8416
$ws4.Range("F20").Value = 293
$ws4.Range("F21").Value = 4
$ws4.Range("F24").Value = 165
$ws4.Range("F41").Value = 15
$ws4.Range("F46").Value = 5185
